$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price (D) and Volume(1h) (E) columns are plain, unformatted text cells
# in the source sheet (t="inlineStr", no numFmt). Re-assigning a bare
# numeric-looking string such as "113.28" through COM .Value would be
# auto-coerced to a real number (and "35.80"/"3.30" would even lose their
# trailing zero), so those are entered with a leading apostrophe - exactly
# like typing them into Excel - to keep them as literal text. Percent
# cells already carry the two-space left/right padding that keeps them
# text in the source file, so they round-trip as-is.

$ws.Range("D2").Value = "49.775.23"
$ws.Range("E2").Value = "  +3.41%  "

$ws.Range("D3").Value = "2.641.43"
$ws.Range("E3").Value = "  +5.59%  "

$ws.Range("E4").Value = "  +0.08%  "

$ws.Range("D5").Value = "'113.28"
$ws.Range("E5").Value = "  +6.42%  "

$ws.Range("D6").Value = "'326.43"
$ws.Range("E6").Value = "  +2.04%  "

$ws.Range("D7").Value = "'0.528"
$ws.Range("E7").Value = "  +0.95%  "

$ws.Range("E8").Value = "  +0.04%  "

$ws.Range("D9").Value = "'0.551"
$ws.Range("E9").Value = "  +2.10%  "

$ws.Range("D10").Value = "'40.74"
$ws.Range("E10").Value = "  +4.89%  "

$ws.Range("E11").Value = "  +0.43%  "

$ws.Range("E12").Value = "  +1.63%  "

$ws.Range("E13").Value = "  +1.02%  "

$ws.Range("D14").Value = "'7.31"
$ws.Range("E14").Value = "  +3.22%  "

$ws.Range("D15").Value = "3.056.05"
$ws.Range("E15").Value = "  +5.63%  "

$ws.Range("D16").Value = "2.642.26"
$ws.Range("E16").Value = "  +5.51%  "

$ws.Range("D17").Value = "'0.867"
$ws.Range("E17").Value = "  +3.97%  "

$ws.Range("D18").Value = "49.688.61"
$ws.Range("E18").Value = "  +3.45%  "

$ws.Range("D19").Value = "'13.09"
$ws.Range("E19").Value = "  +0.75%  "

$ws.Range("D20").Value = "'6.72"
$ws.Range("E20").Value = "  +1.12%  "

$ws.Range("E21").Value = "  -1.89%  "

$ws.Range("D22").Value = "0.0₃0951"
$ws.Range("E22").Value = "  +1.59%  "

$ws.Range("D23").Value = "'72.02"
$ws.Range("E23").Value = "  +1.09%  "

$ws.Range("D24").Value = "'276.20"
$ws.Range("E24").Value = "  +1.39%  "

$ws.Range("E25").Value = "  +2.01%  "

$ws.Range("D26").Value = "'26.63"
$ws.Range("E26").Value = "  +3.09%  "

$ws.Range("E27").Value = "  -0.03%  "

$ws.Range("D28").Value = "'9.95"
$ws.Range("E28").Value = "  +2.17%  "

$ws.Range("E29").Value = "  -3.27%  "

$ws.Range("D30").Value = "'35.80"
$ws.Range("E30").Value = "  +2.20%  "

$ws.Range("E31").Value = "  -1.77%  "

$ws.Range("D32").Value = "'50.33"
$ws.Range("E32").Value = "  +2.02%  "

$ws.Range("D33").Value = "'19.43"
$ws.Range("E33").Value = "  +1.57%  "

$ws.Range("D34").Value = "'5.39"
$ws.Range("E34").Value = "  +1.75%  "

$ws.Range("D35").Value = "'0.0802"
$ws.Range("E35").Value = "  +3.16%  "

$ws.Range("E36").Value = "  -0.21%  "

$ws.Range("D38").Value = "'4.76"
$ws.Range("E38").Value = "  +3.73%  "

$ws.Range("E39").Value = "  +6.01%  "

$ws.Range("D40").Value = "'124.85"
$ws.Range("E40").Value = "  +2.58%  "

$ws.Range("E41").Value = "  +1.34%  "

$ws.Range("E42").Value = "  +1.50%  "

$ws.Range("D43").Value = "'21.84"
$ws.Range("E43").Value = "  -0.90%  "

$ws.Range("D44").Value = "'0.0312"
$ws.Range("E44").Value = "  +2.03%  "

$ws.Range("D45").Value = "2.064.28"
$ws.Range("E45").Value = "  +2.76%  "

$ws.Range("D46").Value = "'3.30"
$ws.Range("E46").Value = "  +3.78%  "

$ws.Range("D47").Value = "'2.28"
$ws.Range("E47").Value = "  +13.95%  "

$ws.Range("E48").Value = "  +4.23%  "

$ws.Range("D49").Value = "'9.06"
$ws.Range("E49").Value = "  +1.31%  "

$ws.Range("E50").Value = "  +3.50%  "

$ws.Range("D51").Value = "'59.18"
$ws.Range("E51").Value = "  +4.24%  "
